# Add three new rows of study-chair data to the "StudyChairs" worksheet.
# The workbook has three sheets: BookShelves, submenuItems, StudyChairs.
# StudyChairs (sheet index 3) currently only has a header row (row 1:
# Model Name | Brand Name | Price). We append rows 2-4 with product data,
# matching the order of cell writes so the shared-string table is built up
# in the same sequence as the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Row 2 (brand + price first, model name filled in last to mirror the
# original authoring order captured in the shared-string table).
$ws.Range("B2").Value = "By Urban Ladder"
$ws.Range("C2").Value = "₹12,287"

# Row 3
$ws.Range("A3").Value = "Galen Study Chair In Black Colour"
$ws.Range("C3").Value = "₹7,505"

# Row 4
$ws.Range("A4").Value = "Hawley Study Chair"
$ws.Range("C4").Value = "₹6,440"

# Row 2 model name filled in after rows 3-4, then the repeated brand name
# for rows 3-4.
$ws.Range("A2").Value = "Mika Leatherette Study Chair In Scarlet Red Colour"
$ws.Range("B3").Value = "By Urban Ladder"
$ws.Range("B4").Value = "By Urban Ladder"

# Leave the active cell on A2, matching the saved selection in the sheet.
$ws.Range("A2").Select()
